$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.412.45'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.939.07'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7482'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '245.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.56%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.53'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3160'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06963'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7795'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07998'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.937.44'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.356'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.52'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.432.75'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '252.52'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007896'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.745'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.190.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.005'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.005'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.673'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.501'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.47'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1322'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('E29').Value = '  -4.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.363'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.513'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.351'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.092'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05149'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.272'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7441'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.781'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01944'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.805'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.418'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4456'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.962'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.005'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.774'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.439'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '979.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +10.44%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.15'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06025'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.50%  '
